$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (player name, position(s), team) for rows 2-19
$names = @(
    "De'Aaron Fox",
    "Miles Bridges",
    "DeMar DeRozan",
    "Goga Bitadze",
    "Evan Mobley",
    "Brook Lopez",
    "Ausar Thompson",
    "Naz Reid",
    "Nikola Vucevic",
    "Josh Giddey",
    "Tyler Herro",
    "Luke Kennard",
    "Mikal Bridges",
    "Onyeka Okongwu",
    "Scottie Barnes",
    "Santi Aldama",
    "Luka Doncic",
    "Ja Morant"
)

$positions = @(
    "PG",
    "PG",
    "SF,PF",
    "C",
    "PF,C",
    "C",
    "SF,PF",
    "PF,C",
    "PF,C",
    "PG,SG,SF",
    "PG,SG",
    "SG",
    "SG,SF,PF",
    "PF,C",
    "PG,SG,SF,PF",
    "PF,C",
    "PG,SG",
    "PG"
)

$teams = @(
    "Sacramento Kings",
    "Charlotte Hornets",
    "Sacramento Kings",
    "Orlando Magic",
    "Cleveland Cavaliers",
    "Milwaukee Bucks",
    "Detroit Pistons",
    "Minnesota Timberwolves",
    "Chicago Bulls",
    "Chicago Bulls",
    "Miami Heat",
    "Memphis Grizzlies",
    "New York Knicks",
    "Atlanta Hawks",
    "Toronto Raptors",
    "Memphis Grizzlies",
    "Dallas Mavericks",
    "Memphis Grizzlies"
)

# Step 1: clear out the old data cells (the header row in A1:C1 - and its
# bold/bordered style - is left untouched) so none of the stale values linger
# if the new list were ever shorter than the old one.
$ws.Range("A2:C19").ClearContents()

# Step 2: write the real values column-by-column (all of column A, then all
# of column B, then all of column C) so that values are populated in the
# same grouping the source roster data came in.
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value2 = $names[$i]
}
for ($i = 0; $i -lt $positions.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value2 = $positions[$i]
}
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value2 = $teams[$i]
}
